# GateMaster.pptx edit:
# Insert a new "Frontend Snapshots" slide right before the existing
# "ER Diagram and Tables" slide (which becomes the new last slide).

$p = $ppt.ActivePresentation

# The "ER Diagram and Tables" slide is currently the last slide (index 7).
# Insert the new slide at that position, using the same "Title and Content"
# style layout that slide uses (ppLayoutText = 2), which pushes the ER
# Diagram slide down to become the new last slide.
$erIndex = $p.Slides.Count
$newSlide = $p.Slides.Add($erIndex, 2)

# Title placeholder: "Frontend Snapshots" - bold, centered, same custom
# font used by the sibling section-title slides in this deck.
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Frontend Snapshots"
$title.Font.Bold = 1
$title.Font.Name = "Footlight MT Light"
$title.ParagraphFormat.Alignment = 2

# Content placeholder: one bullet paragraph per frontend page/section.
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Home`rRegister / Login`rStreams`rAbout Gate`rFeedback Review"
